{"js": "// \"fix save file autofill\" \u2014 the payer's name gets auto-filled in, while\n// the other fields (which used to hold sample/default data) are reset back\n// to their blank \"Label: \" placeholder text.\nconst body = context.document.body;\n\nasync function findFirst(searchText) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n  return results.items[0];\n}\n\n// 1) Auto-fill the requester's name after the existing label. Appending\n//    (rather than replacing the whole run) keeps the run's xml:space\n//    handling consistent with a simple fill-in.\nconst hoTenRange = await findFirst(\"H\u1ecd, t\u00ean ng\u01b0\u1eddi \u0111\u1ec1 ngh\u1ecb thanh to\u00e1n: \");\nhoTenRange.insertText(\"\u0110\u1ed7 Thanh B\u00ecnh\", Word.InsertLocation.end);\nawait context.sync();\n\n// 2) Clear \"\u0110\u01a1n v\u1ecb\" back to just the label.\nconst donViRange = await findFirst(\"\u0110\u01a1n v\u1ecb: Khoa C\u00f4ng Ngh\u1ec7 Th\u00f4ng Tin\");\ndonViRange.insertText(\"\u0110\u01a1n v\u1ecb: \", Word.InsertLocation.replace);\nawait context.sync();\n\n// 3) Clear \"N\u1ed9i dung thanh to\u00e1n\" back to just the label.\nconst noiDungRange = await findFirst(\n  \"N\u1ed9i dung thanh to\u00e1n: Thanh to\u00e1n kinh ph\u00ed c\u1ee7a \u0111\u1ec3 t\u00e0i c\u1ea5p c\u01a1 s\u1edf Nghi\u00ean c\u1ee9u ki\u1ebfn tr\u00fac chip x\u1eed l\u00fd m\u1eadt m\u00e3 theo ti\u00eau chu\u1ea9n Trusted Platform Module 2.0 (TPM 2.0) c\u1ee7a Trusted Computing Group (TCG)\\\" M\u00e3 s\u1ed1 \u0111\u1ec3 t\u00e0i 19/2023/CS do L\u00ea Anh Ti\u1ebfn l\u00e0m ch\u1ee7 nhi\u1ec7m.\"\n);\nnoiDungRange.insertText(\"N\u1ed9i dung thanh to\u00e1n: \", Word.InsertLocation.replace);\nawait context.sync();\n\n// 4) Clear \"M\u00e3 s\u1ed1 \u0111\u1ec1 t\u00e0i\" back to just the label.\nconst maSoRange = await findFirst(\"M\u00e3 s\u1ed1 \u0111\u1ec1 t\u00e0i: 19/2023/CS do L\u00ea Anh Ti\u1ebfn l\u00e0m ch\u1ee7 nhi\u1ec7m.\");\nmaSoRange.insertText(\"M\u00e3 s\u1ed1 \u0111\u1ec1 t\u00e0i: \", Word.InsertLocation.replace);\nawait context.sync();\n\n// 5) Clear \"S\u1ed1 ti\u1ec1n\" back to just the label.\nconst soTienRange = await findFirst(\"S\u1ed1 ti\u1ec1n: 50.000.000 \u0111\u1ed3ng\");\nsoTienRange.insertText(\"S\u1ed1 ti\u1ec1n: \", Word.InsertLocation.replace);\nawait context.sync();\n\n// 6) Clear \"Vi\u1ebft b\u1eb1ng ch\u1eef\" back to just the label.\nconst vietBangChuRange = await findFirst(\"Vi\u1ebft b\u1eb1ng ch\u1eef: N\u0103m m\u01b0\u01a1i tri\u1ec7u \u0111\u1ed3ng ch\u1eb3n.\");\nvietBangChuRange.insertText(\"Vi\u1ebft b\u1eb1ng ch\u1eef: \", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# \"fix save file autofill\" \u2014 the payer's name gets auto-filled in, while\n# the other fields (which used to hold sample/default data) are reset back\n# to their blank \"Label: \" placeholder text.\n\n$d = $word.ActiveDocument\n\nfunction Replace-DocText {\n    param(\n        [string]$SearchText,\n        [string]$ReplaceText\n    )\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    $found = $find.Execute(\n        $SearchText,   # FindText\n        $false,        # MatchCase\n        $false,        # MatchWholeWord\n        $false,        # MatchWildcards\n        $false,        # MatchSoundsLike\n        $false,        # MatchAllWordForms\n        $true,         # Forward\n        1,             # Wrap (wdFindContinue)\n        $false,        # Format\n        $ReplaceText,  # ReplaceWith\n        2              # Replace (wdReplaceOne)\n    )\n\n    if (-not $found) {\n        throw \"Text not found: $SearchText\"\n    }\n}\n\n# 1) Auto-fill the requester's name after the existing label.\nReplace-DocText \"H\u1ecd, t\u00ean ng\u01b0\u1eddi \u0111\u1ec1 ngh\u1ecb thanh to\u00e1n: \" \"H\u1ecd, t\u00ean ng\u01b0\u1eddi \u0111\u1ec1 ngh\u1ecb thanh to\u00e1n: \u0110\u1ed7 Thanh B\u00ecnh\"\n\n# 2) Clear \"\u0110\u01a1n v\u1ecb\" back to just the label.\nReplace-DocText \"\u0110\u01a1n v\u1ecb: Khoa C\u00f4ng Ngh\u1ec7 Th\u00f4ng Tin\" \"\u0110\u01a1n v\u1ecb: \"\n\n# 3) Clear \"N\u1ed9i dung thanh to\u00e1n\" back to just the label.\nReplace-DocText \"N\u1ed9i dung thanh to\u00e1n: Thanh to\u00e1n kinh ph\u00ed c\u1ee7a \u0111\u1ec3 t\u00e0i c\u1ea5p c\u01a1 s\u1edf Nghi\u00ean c\u1ee9u ki\u1ebfn tr\u00fac chip x\u1eed l\u00fd m\u1eadt m\u00e3 theo ti\u00eau chu\u1ea9n Trusted Platform Module 2.0 (TPM 2.0) c\u1ee7a Trusted Computing Group (TCG)`\" M\u00e3 s\u1ed1 \u0111\u1ec3 t\u00e0i 19/2023/CS do L\u00ea Anh Ti\u1ebfn l\u00e0m ch\u1ee7 nhi\u1ec7m.\" \"N\u1ed9i dung thanh to\u00e1n: \"\n\n# 4) Clear \"M\u00e3 s\u1ed1 \u0111\u1ec1 t\u00e0i\" back to just the label.\nReplace-DocText \"M\u00e3 s\u1ed1 \u0111\u1ec1 t\u00e0i: 19/2023/CS do L\u00ea Anh Ti\u1ebfn l\u00e0m ch\u1ee7 nhi\u1ec7m.\" \"M\u00e3 s\u1ed1 \u0111\u1ec1 t\u00e0i: \"\n\n# 5) Clear \"S\u1ed1 ti\u1ec1n\" back to just the label.\nReplace-DocText \"S\u1ed1 ti\u1ec1n: 50.000.000 \u0111\u1ed3ng\" \"S\u1ed1 ti\u1ec1n: \"\n\n# 6) Clear \"Vi\u1ebft b\u1eb1ng ch\u1eef\" back to just the label.\nReplace-DocText \"Vi\u1ebft b\u1eb1ng ch\u1eef: N\u0103m m\u01b0\u01a1i tri\u1ec7u \u0111\u1ed3ng ch\u1eb3n.\" \"Vi\u1ebft b\u1eb1ng ch\u1eef: \"\n"}
